# Add two new columns, I (I0) and J (IF), to the save-data worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) -------------------------------------------------
# Copy the formatting of the existing "IP" header (H1) onto the two new
# header cells so they match the bold/centered header style, then set
# their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# --- Data rows (rows 2-36) --------------------------------------------
# For every data row, column J mirrors the existing "IP" value in
# column H, and column I is 1 -- except for the most recent game
# (row 2), where I0=7 and IF=9.
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($r, 8).Value2
}
$ws.Cells.Item(2, 9).Value2 = 7
$ws.Cells.Item(2, 10).Value2 = 9
